$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 44: TP 9 / TP 10 entries for Oregon tide pools (reuse formatting from row 43)
[void]$ws.Range("A43:F43").Copy()
[void]$ws.Range("A44").PasteSpecial(-4122)   # xlPasteFormats - carry over date style, etc.

$ws.Range("A44").Value = 43715
$ws.Range("B44").Value = 2212.9281273883298
$ws.Range("C44").Value = 2207.0300000000002
$ws.Range("D44").Formula = "=100*(B44-C44)/C44"
$ws.Range("E44").Value = 169
$ws.Range("F44").Value = "Opened CRM (8/30/2019)"

# Keep the current selection pointed at F43, matching Excel's behavior after
# entering the new row of data (selection was F42:F43 before the edit).
[void]$ws.Range("F43").Select()
